$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the three new "record" header columns (AD:AF), copying the existing
# header style (bold, bordered, centered) from the adjacent header cell.
$ws.Range("AC1").Copy($ws.Range("AD1:AF1"))
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Every player row (2-48) shares the same team season record.
for ($r = 2; $r -le 48; $r++) {
    $ws.Cells.Item($r, 30).Value = 74
    $ws.Cells.Item($r, 31).Value = 88
    $ws.Cells.Item($r, 32).Value = 0
}
